$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a pre-existing data bug: I5 should read "facilities" (it duplicated the
# "20" room-count value by mistake).
$ws.Range("I5").Value = "facilities"

# Row 6 was filler ("test") - replace with real hotel #105 data.
# Keep it formatted as Text like the rest of the table (column style).
$ws.Range("A6:I6").NumberFormat = "@"
$ws.Range("A6").Value = "105"
$ws.Range("B6").Value = "The Wyndham New Yorker Hotel"
$ws.Range("C6").Value = "United States"
$ws.Range("D6").Value = "New York"
$ws.Range("E6").Value = "481 Eighth Avenue"
$ws.Range("F6").Value = "8888881"
$ws.Range("G6").Value = "5"
$ws.Range("H6").Value = "300"
$ws.Range("I6").Value = "Fitness"

# Row 7 was filler ("ttt") - replace with real hotel #106 data.
$ws.Range("A7").Value = "106"
$ws.Range("B7").Value = "Hotel Mela Times Square"
$ws.Range("C7").Value = "United States"
$ws.Range("D7").Value = "New York"
$ws.Range("E7").Value = "120 West 44th Street"
$ws.Range("F7").Value = "1829981"
$ws.Range("G7").Value = "5"
$ws.Range("H7").Value = "600"
$ws.Range("I7").Value = "everything"

# Row 8 was filler ("dddddd") - replace with real hotel #107 data.
$ws.Range("A8").Value = "107"
$ws.Range("B8").Value = "Hotel Pennsylvania"
$ws.Range("C8").Value = "United States"
$ws.Range("D8").Value = "New York"
$ws.Range("E8").Value = "401 7th Avenue"
$ws.Range("F8").Value = "999999"
$ws.Range("G8").Value = "5"
$ws.Range("H8").Value = "400"
$ws.Range("I8").Value = "everything"

# Row 9 was filler ("ppppp") - replace with real hotel #108 data.
$ws.Range("A9").Value = "108"
$ws.Range("B9").Value = "The Savoy Hotel"
$ws.Range("C9").Value = "United Kingdom"
$ws.Range("D9").Value = "London"
$ws.Range("E9").Value = "Strand, West End Soho"
$ws.Range("F9").Value = "11233"
$ws.Range("G9").Value = "3"
$ws.Range("H9").Value = "40"
$ws.Range("I9").Value = "---"

# New row 10 - hotel #109.
$ws.Range("A10").Value = "109"
$ws.Range("B10").Value = "Days Inn Hilton Head"
$ws.Range("C10").Value = "United States"
$ws.Range("D10").Value = "hilton"
$ws.Range("E10").Value = "9 Marina Side Drive"
$ws.Range("F10").Value = "999999"
$ws.Range("G10").Value = "5"
$ws.Range("H10").Value = "200"
$ws.Range("I10").Value = "fitness"

# New row 11 - hotel #110.
$ws.Range("A11").Value = "110"
$ws.Range("B11").Value = "Hilton Head Marriott Resort & Spa"
$ws.Range("C11").Value = "United States"
$ws.Range("D11").Value = "Hilton"
$ws.Range("E11").Value = "---"
$ws.Range("F11").Value = "11111"
$ws.Range("G11").Value = "5"
$ws.Range("H11").Value = "500"
$ws.Range("I11").Value = "unknown"

# Leave the selection on the last data cell entered, matching the author's
# final cursor position.
$ws.Range("I10").Select()
